{"js": "// Locate the existing bullet paragraph under \"Siege Analytics\" that we need\n// to insert the three new bullet points after.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Lead comprehensive research initiatives\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nconst newLines = [\n  \"\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\",\n  \"\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\",\n  \"\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations $5M+ and enabling smaller nonprofits to conduct redistricting analysis\"\n];\n\nlet anchor = target;\nfor (const line of newLines) {\n  anchor = anchor.insertParagraph(line, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the existing bullet paragraph under \"Siege Analytics\" that we need\n# to insert the three new bullet points after.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Lead comprehensive research initiatives*\") {\n        $target = $p\n        break\n    }\n}\n\n$newLines = @(\n    \"\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\",\n    \"\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\",\n    \"\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis\"\n)\n\n$idx = $target.Index\nforeach ($line in $newLines) {\n    $d.Paragraphs($idx).Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs($idx + 1)\n    $newPara.Range.Text = $line\n    $idx = $idx + 1\n}\n"}
